# =====================================================================
# Edit script: "The Woven Tapestry of Evolution" -> "Navigating the
# Enigmatic Realm of Mathematics" rewrite, per commit diff.
# =====================================================================

$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceAll = 2
$LB = [char]11   # Word "manual line break" (serializes to <w:br/>)

function ReplaceText($old, $new) {
    # Find the (unique) old text and overwrite the Range's .Text
    # directly, rather than going through Find.Execute's own
    # Replacement text path -- the latter silently runs "smart quotes"
    # autocorrect over the replacement string (straight "'" -> U+2019),
    # which would corrupt text like "nature's". Assigning .Text keeps
    # the run's formatting (rPr) and inserts the literal characters.
    $r = $d.Content
    $r.Find.Execute($old, $true) | Out-Null
    $r.Text = $new
}

function InsertAfterAnchor($anchorText, $newText) {
    # Locate the (unique) anchor text, collapse the found range to its
    # end, and insert newText there so it inherits the formatting of
    # the text it is being appended after (matching real Word typing
    # behaviour at that caret position).
    $r = $d.Content
    $r.Find.Execute($anchorText, $true) | Out-Null
    $r.Collapse(0) | Out-Null
    $r.InsertAfter($newText) | Out-Null
}

# ---------------------------------------------------------------------
# Title / byline / e-mail
# ---------------------------------------------------------------------
ReplaceText "The Woven Tapestry of Evolution" "Navigating the Enigmatic Realm of Mathematics"
ReplaceText " Allyson Hill" " Emily Carter"
ReplaceText "allysonhill@biosphere" "emilycarter@columbiaacademy"
ReplaceText "institute" "org"

# ---------------------------------------------------------------------
# Body paragraph (Introduction block)
# ---------------------------------------------------------------------
ReplaceText "Across the expansive landscape of sciences, from the intricate workings of biology to the cosmic tapestry of astronomy, lies the undeniable truth of evolution" `
            "Welcome to the captivating realm of Mathematics! This realm holds immense power, beauty, and mesmerizing patterns that shape our understanding of the world"

ReplaceText " Evolution is not merely a concept or a theory; it is a compelling narrative etched into the very fabric of life, a testament to the dynamic nature of existence" `
            " Mathematics is a profound discipline that not only challenges our intellect but also serves as a tool for unlocking a plethora of knowledge domains"

InsertAfterAnchor "unlocking a plethora of knowledge domains" (
  "." + " In this essay, we shall embark on an enlightening voyage through Mathematics, unveiling its indelible imprint on various spheres of life, its intricate beauty, and the Limitless possibilities it offers")

ReplaceText "In the realm of biology, evolution unveils the magnificent diversity of life, from the single-celled organisms that grace the microscopic world to the intricate web of ecosystems that span continents" `
            "Mathematics weaves its intricate web throughout the very fabric of our existence, from the rhythms of nature's designs to the symphony of numbers reverberating within our universe"

ReplaceText " Each species bears the unique imprint of its evolutionary journey, a testament to adaptation and resilience in the face of ever-changing environments" `
            " It provides the solid foundation upon which the pillars of physics, engineering, and astronomy stand tall, enabling us to explore the enigmatic wonders of the cosmos"

InsertAfterAnchor "explore the enigmatic wonders of the cosmos" (
  "." + " Moreover, it forms the unrelenting backbone of economics, finance, and architecture, empowering us to comprehend the intricate complexities of our world")

ReplaceText "Beyond biology, evolution extends its influence to fields such as geology and paleontology" `
            "Mathematics is not merely a collection of abstract theorems and equations; it is a profound language that unveils the symphony of patterns concealed within seemingly unrelated phenomena"

ReplaceText " Sedimentary layers, like the rings of a tree trunk, hold the encoded secrets of ancient climates and the transformation of landscapes over time" `
            " The intricate dance of numbers and shapes captivates our imagination, inspiring awe and wonder, etching an enduring tapestry of intellectual beauty that sets us apart from other beings"

ReplaceText " Fossils, preserved relics of past life, offer tantalizing glimpses into vanished worlds and extinct species, painting a vivid picture of Earth's dynamic history" `
            " It is the enigma that entices seekers of knowledge, leading them down a path of discovery that stretches beyond the boundaries of time"

# New "Body:" sub-heading + two paragraphs of new material, inserted
# right after the sentence above (still inside the same w:p).
InsertAfterAnchor "leading them down a path of discovery that stretches beyond the boundaries of time" (
  "." + $LB + $LB + "Body:" + $LB + $LB +
  "Mathematics is the language that unlocks the secrets of the universe, from the tiniest atoms to the vast expanse of constellations that twinkle above" + "." +
  " Its myriad applications span a kaleidoscope of fields, including engineering, where it ensures the safety and integrity of structures, from sky-piercing skyscrapers to the grace of bridges arching over rivers" + "." +
  " It also finds its way into economics, where it unveils patterns and trends in financial markets, empowering us to make informed decisions and navigate the intricacies of commerce" + "." + $LB + $LB +
  "Mathematics is an integral part of the scientific fabric that unveils the mysteries of our universe" + "." +
  " It is the Rosetta stone of science, enabling scientists to translate and decipher the concealed language of nature" + "." +
  " Through mathematical equations, we probe the depths of quantum mechanics, unraveling the enigmatic behavior of particles and forces" + "." +
  " We delve into the intricacies of biology, discovering the inner workings of cells and the symphony of life's mechanisms" + "." +
  " With every new equation, mathematics serves as a torch illuminating our understanding of the world's boundless wonders" + "." + $LB + $LB +
  "Moreover, mathematics possesses an inherent aesthetic allure that transcends its practical applications" + "." +
  " Its elegant patterns and intricate symmetries dance before our eyes, sparking a profound sense of wonder and appreciation" + "." +
  " The interplay of numbers, shapes, and patterns evokes an artistic resonance that captivates the hearts and minds of mathematicians and non-mathematicians alike" + "." +
  " It reminds us that Mathematics is not merely a set of rigid rules and formulas but also a canvas of boundless creativity and imagination")

# ---------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------
ReplaceText "Evolution stands as a unifying force across diverse scientific disciplines, showcasing the intricate interconnectedness of life and the cosmos" `
            "Mathematics is a realm of boundless exploration, beauty, and limitless potential"

ReplaceText " It reminds us of the constant flux of change, the tapestry of life ever-weaving its exquisite patterns, and the awe-inspiring journey we are all part of" `
            " It shapes our understanding of the world, unraveling the enigmatic patterns that define our existence"

InsertAfterAnchor "unraveling the enigmatic patterns that define our existence" (
  "." +
  " Its applications span the vast spectrum of human endeavors, ranging from science and engineering to economics and the arts" + "." +
  " Mathematics ignites our intellect, kindles our curiosity, and unveils the wonders of the universe" + "." +
  " It challenges us to push the boundaries of knowledge and innovation, paving the way for transformative discoveries that shape our future" + "." +
  " Through Mathematics, humanity can unlock its full potential and illuminate the path toward progress, prosperity, and fulfillment")

# ---------------------------------------------------------------------
# Trailing empty paragraph added at the end of the document body.
# ---------------------------------------------------------------------
$d.Content.InsertParagraphAfter() | Out-Null
